$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2704.2856
$ws.Range("I17").Value = 3813
$ws.Range("J17").Value = 1872.75
$ws.Range("K17").Value = 11439
$ws.Range("L17").Value = 5618.25
$ws.Range("M17").Value = -11271
$ws.Range("N17").Value = -5954.25
$ws.Range("H43").Value = 4920.8335
$ws.Range("I43").Value = 1750
$ws.Range("J43").Value = 8091.6665
$ws.Range("K43").Value = 1750
$ws.Range("L43").Value = 8091.6665
$ws.Range("M43").Value = -1681
$ws.Range("N43").Value = -8229.666499999999
$ws.Range("H61").Value = 750
$ws.Range("I61").Value = 750
$ws.Range("K61").Value = 2250
$ws.Range("M61").Value = -2078
$ws.Range("H62").Value = 6183.5386
$ws.Range("I62").Value = 4064.3333
$ws.Range("K62").Value = 4064.3333
$ws.Range("M62").Value = -3440.3333
$ws.Range("H65").Value = 6183.5386
$ws.Range("I65").Value = 4064.3333
$ws.Range("K65").Value = 20321.6665
$ws.Range("M65").Value = -17201.6665
$ws.Range("H70").Value = 92551.22
$ws.Range("I70").Value = 1329.6666
$ws.Range("J70").Value = 138162
$ws.Range("K70").Value = 3988.9998
$ws.Range("L70").Value = 414486
$ws.Range("M70").Value = -3718.9998
$ws.Range("N70").Value = -415026
$ws.Range("H73").Value = 92551.22
$ws.Range("I73").Value = 1329.6666
$ws.Range("J73").Value = 138162
$ws.Range("K73").Value = 3988.9998
$ws.Range("L73").Value = 414486
$ws.Range("M73").Value = -3052.9998
$ws.Range("N73").Value = -416358
$ws.Range("H106").Value = 25983.777
$ws.Range("I106").Value = 32494
$ws.Range("K106").Value = 32494
$ws.Range("M106").Value = -31863
$ws.Range("H116").Value = 7180.1
$ws.Range("I116").Value = 6183.6665
$ws.Range("K116").Value = 6183.6665
$ws.Range("M116").Value = -2741.6665
$ws.Range("H125").Value = 927.8
$ws.Range("I125").Value = 927.8
$ws.Range("K125").Value = 8350.199999999999
$ws.Range("M125").Value = -5890.199999999999
$ws.Range("H132").Value = 1317.7142
$ws.Range("I132").Value = 1256.6666
$ws.Range("K132").Value = 3769.9998
$ws.Range("M132").Value = -1239.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15348.5
$ws.Range("I32").Value = 6403.727
$ws.Range("K32").Value = 6403.727
$ws.Range("M32").Value = -6116.727
$ws.Range("H74").Value = 1564.1
$ws.Range("I74").Value = 750.3333
$ws.Range("K74").Value = 750.3333
$ws.Range("M74").Value = 123.6667
$ws.Range("H76").Value = 73300
$ws.Range("J76").Value = 73300
$ws.Range("L76").Value = 73300
$ws.Range("N76").Value = -73976
$ws.Range("H77").Value = 1564.1
$ws.Range("I77").Value = 750.3333
$ws.Range("K77").Value = 3751.6665
$ws.Range("M77").Value = 616.3334999999997
$ws.Range("H79").Value = 73300
$ws.Range("J79").Value = 73300
$ws.Range("L79").Value = 73300
$ws.Range("N79").Value = -75640
$ws.Range("H97").Value = 1002.25
$ws.Range("I97").Value = 336.33334
$ws.Range("K97").Value = 336.33334
$ws.Range("M97").Value = 159.66666
$ws.Range("H122").Value = 373010.66
$ws.Range("I122").Value = 557046.9
$ws.Range("K122").Value = 1671140.7
$ws.Range("M122").Value = -1668690.7
$ws.Range("H132").Value = 2129.375
$ws.Range("I132").Value = 1938
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5814
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3284
$ws.Range("N132").Value = -20060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3332.1924
$ws.Range("J99").Value = 4498.3335
$ws.Range("L99").Value = 4498.3335
$ws.Range("N99").Value = -7494.3335
$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1480
$ws.Range("J16").Value = 2999
$ws.Range("L16").Value = 2999
$ws.Range("N16").Value = -3573
$ws.Range("H62").Value = 59714
$ws.Range("I62").Value = 2999.5
$ws.Range("K62").Value = 2999.5
$ws.Range("M62").Value = -2375.5
$ws.Range("H65").Value = 59714
$ws.Range("I65").Value = 2999.5
$ws.Range("K65").Value = 14997.5
$ws.Range("M65").Value = -11877.5
$ws.Range("H99").Value = 17837.4
$ws.Range("I99").Value = 15748.25
$ws.Range("J99").Value = 18597.092
$ws.Range("K99").Value = 15748.25
$ws.Range("L99").Value = 18597.092
$ws.Range("M99").Value = -14250.25
$ws.Range("N99").Value = -21593.092
$ws.Range("H105").Value = 2140.5
$ws.Range("I105").Value = 981
$ws.Range("J105").Value = 3300
$ws.Range("K105").Value = 981
$ws.Range("L105").Value = 3300
$ws.Range("M105").Value = 766
$ws.Range("N105").Value = -6794
$ws.Range("H107").Value = 1241.1428
$ws.Range("I107").Value = 948
$ws.Range("K107").Value = 948
$ws.Range("M107").Value = 972
$ws.Range("H113").Value = 1480
$ws.Range("J113").Value = 2999
$ws.Range("L113").Value = 2999
$ws.Range("N113").Value = -7339
$ws.Range("H126").Value = 17837.4
$ws.Range("I126").Value = 15748.25
$ws.Range("J126").Value = 18597.092
$ws.Range("K126").Value = 47244.75
$ws.Range("L126").Value = 55791.276
$ws.Range("M126").Value = -44774.75
$ws.Range("N126").Value = -60731.276

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3487.75
$ws.Range("J131").Value = 9454.5
$ws.Range("L131").Value = 28363.5
$ws.Range("N131").Value = -38443.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H126").Value = 4496.125
$ws.Range("I126").Value = 3085
$ws.Range("J126").Value = 4966.5
$ws.Range("K126").Value = 9255
$ws.Range("L126").Value = 14899.5
$ws.Range("M126").Value = -6785
$ws.Range("N126").Value = -19839.5
$ws.Range("H132").Value = 4113.5
$ws.Range("I132").Value = 2576.6667
$ws.Range("J132").Value = 6418.75
$ws.Range("K132").Value = 7730.000100000001
$ws.Range("L132").Value = 19256.25
$ws.Range("M132").Value = -5200.000100000001
$ws.Range("N132").Value = -24316.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3899.6667
$ws.Range("I16").Value = 3899.6667
$ws.Range("K16").Value = 3899.6667
$ws.Range("M16").Value = -3729.6667
$ws.Range("H46").Value = 3042.0715
$ws.Range("J46").Value = 3557
$ws.Range("L46").Value = 3557
$ws.Range("N46").Value = -3933
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 3876.0833
$ws.Range("I132").Value = 3610.682
$ws.Range("J132").Value = 4293.143
$ws.Range("K132").Value = 10832.046
$ws.Range("L132").Value = 12879.429
$ws.Range("M132").Value = -8302.045999999998
$ws.Range("N132").Value = -17939.429

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 70000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 70000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H126").Value = 2276.1
$ws.Range("I126").Value = 862.3333
$ws.Range("K126").Value = 2586.9999
$ws.Range("M126").Value = -116.9998999999998
$ws.Range("H132").Value = 1513.421
$ws.Range("I132").Value = 1347.5
$ws.Range("K132").Value = 4042.5
$ws.Range("M132").Value = -1512.5
